$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 405-422: values were re-shuffled/changed ---
# (two new records inserted at the top of this date-block; D/J/K/L/M/N/P shift down)

$ws.Range("D405").Value = 45075
$ws.Range("J405").Value = 50
$ws.Range("K405").Value = 17000
$ws.Range("L405").Value = 17000
$ws.Range("M405").Value = 17000
$ws.Range("P405").Value = 1700

$ws.Range("D406").Value = 45075
$ws.Range("K406").Value = 18000
$ws.Range("L406").Value = 18000
$ws.Range("M406").Value = 18000
$ws.Range("N406").Value = '$/malla 10 kilos'
$ws.Range("P406").Value = 1800

$ws.Range("D407").Value = 44624
$ws.Range("K407").Value = 20000
$ws.Range("L407").Value = 21000
$ws.Range("M407").Value = 20500
$ws.Range("P407").Value = 2050

$ws.Range("D408").Value = 44218
$ws.Range("J408").Value = 50
$ws.Range("K408").Value = 12500
$ws.Range("L408").Value = 13000
$ws.Range("M408").Value = 12800
$ws.Range("N408").Value = '$/caja 10 kilos'
$ws.Range("P408").Value = 1280

$ws.Range("D409").Value = 44921
$ws.Range("K409").Value = 15000
$ws.Range("L409").Value = 16000
$ws.Range("M409").Value = 15500
$ws.Range("P409").Value = 1550

$ws.Range("D410").Value = 44921
$ws.Range("J410").Value = 60
$ws.Range("K410").Value = 17000
$ws.Range("L410").Value = 18000
$ws.Range("M410").Value = 17500
$ws.Range("P410").Value = 1750

$ws.Range("D411").Value = 45014
$ws.Range("J411").Value = 60
$ws.Range("K411").Value = 17000
$ws.Range("L411").Value = 18000
$ws.Range("M411").Value = 17500
$ws.Range("P411").Value = 1750

$ws.Range("D412").Value = 45014
$ws.Range("J412").Value = 50
$ws.Range("K412").Value = 20000
$ws.Range("L412").Value = 20000
$ws.Range("M412").Value = 20000
$ws.Range("N412").Value = '$/malla 10 kilos'
$ws.Range("P412").Value = 2000

$ws.Range("D413").Value = 44648
$ws.Range("J413").Value = 40
$ws.Range("K413").Value = 21000
$ws.Range("L413").Value = 22000
$ws.Range("M413").Value = 21500
$ws.Range("P413").Value = 2150

$ws.Range("D414").Value = 44469
$ws.Range("J414").Value = 60
$ws.Range("K414").Value = 16000
$ws.Range("L414").Value = 17000
$ws.Range("M414").Value = 16500
$ws.Range("P414").Value = 1650

$ws.Range("D415").Value = 44651
$ws.Range("J415").Value = 60
$ws.Range("K415").Value = 19000
$ws.Range("L415").Value = 20000
$ws.Range("M415").Value = 19500
$ws.Range("N415").Value = '$/caja 10 kilos'
$ws.Range("P415").Value = 1950

$ws.Range("D416").Value = 44160
$ws.Range("J416").Value = 110
$ws.Range("K416").Value = 7500
$ws.Range("L416").Value = 8000
$ws.Range("M416").Value = 7773
$ws.Range("N416").Value = '$/caja 10 kilos'
$ws.Range("P416").Value = 777

$ws.Range("D417").Value = 44988
$ws.Range("J417").Value = 50
$ws.Range("K417").Value = 21000
$ws.Range("L417").Value = 21000
$ws.Range("M417").Value = 21000
$ws.Range("N417").Value = '$/malla 10 kilos'
$ws.Range("P417").Value = 2100

$ws.Range("D418").Value = 44999
$ws.Range("J418").Value = 40
$ws.Range("K418").Value = 20000
$ws.Range("L418").Value = 21000
$ws.Range("M418").Value = 20500
$ws.Range("P418").Value = 2050

$ws.Range("D419").Value = 45040
$ws.Range("J419").Value = 60
$ws.Range("K419").Value = 17000
$ws.Range("L419").Value = 18000
$ws.Range("M419").Value = 17500
$ws.Range("P419").Value = 1750

$ws.Range("D420").Value = 45040
$ws.Range("J420").Value = 60
$ws.Range("K420").Value = 19000
$ws.Range("L420").Value = 20000
$ws.Range("M420").Value = 19500
$ws.Range("N420").Value = '$/malla 10 kilos'
$ws.Range("P420").Value = 1950

$ws.Range("D421").Value = 44662
$ws.Range("J421").Value = 100
$ws.Range("K421").Value = 19000
$ws.Range("L421").Value = 20000
$ws.Range("M421").Value = 19500
$ws.Range("P421").Value = 1950

$ws.Range("D422").Value = 44473
$ws.Range("J422").Value = 40
$ws.Range("K422").Value = 16000
$ws.Range("L422").Value = 17000
$ws.Range("M422").Value = 16500
$ws.Range("N422").Value = '$/caja 10 kilos'
$ws.Range("P422").Value = 1650

# --- Append two new rows (423, 424) at the end of the data block ---

# Row 423
$ws.Range("A423").Value = 7
$ws.Range("B423").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C423").Value = 'Ñuble'
$ws.Range("D423").Value = 44910
$ws.Range("E423").Value = 16
$ws.Range("F423").Value = 100112003
$ws.Range("G423").Value = 'Ajo'
$ws.Range("H423").Value = 'Chino'
$ws.Range("I423").Value = 'Primera'
$ws.Range("J423").Value = 60
$ws.Range("K423").Value = 15000
$ws.Range("L423").Value = 16000
$ws.Range("M423").Value = 15500
$ws.Range("N423").Value = '$/caja 10 kilos'
$ws.Range("O423").Value = 'China'
$ws.Range("P423").Value = 1550
$ws.Range("Q423").Value = 10
$ws.Range("R423").Value = 'Hortaliza'
$ws.Range("D423").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 424
$ws.Range("A424").Value = 7
$ws.Range("B424").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C424").Value = 'Ñuble'
$ws.Range("D424").Value = 44910
$ws.Range("E424").Value = 16
$ws.Range("F424").Value = 100112003
$ws.Range("G424").Value = 'Ajo'
$ws.Range("H424").Value = 'Chino'
$ws.Range("I424").Value = 'Primera'
$ws.Range("J424").Value = 60
$ws.Range("K424").Value = 17000
$ws.Range("L424").Value = 18000
$ws.Range("M424").Value = 17500
$ws.Range("N424").Value = '$/malla 10 kilos'
$ws.Range("O424").Value = 'China'
$ws.Range("P424").Value = 1750
$ws.Range("Q424").Value = 10
$ws.Range("R424").Value = 'Hortaliza'
$ws.Range("D424").NumberFormat = "YYYY-MM-DD HH:MM:SS"
